$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.418.56"
$ws.Range("E2").Value = "  -7.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.688.73"
$ws.Range("E3").Value = "  -5.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.14"
$ws.Range("E5").Value = "  -5.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5087"
$ws.Range("E6").Value = "  -13.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2675"
$ws.Range("E8").Value = "  -3.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.08"
$ws.Range("E9").Value = "  -5.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06319"
$ws.Range("E10").Value = "  -6.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07378"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.689.61"
$ws.Range("E12").Value = "  -5.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.534"
$ws.Range("E13").Value = "  -5.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5793"
$ws.Range("E14").Value = "  -5.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.918.19"
$ws.Range("E15").Value = "  -5.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008665"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.26"
$ws.Range("E17").Value = "  -13.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.474.38"
$ws.Range("E18").Value = "  -7.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.005"
$ws.Range("E19").Value = "  -7.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "186.58"
$ws.Range("E22").Value = "  -10.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.272"
$ws.Range("E23").Value = "  -7.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -5.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.509"
$ws.Range("E26").Value = "  -5.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1171"
$ws.Range("E27").Value = "  -7.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.85"
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05755"
$ws.Range("E30").Value = "  -5.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.329"
$ws.Range("E31").Value = "  -6.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.524"
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.532"
$ws.Range("E33").Value = "  -6.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.665"
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.015"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5964"
$ws.Range("E36").Value = "  -6.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.355"
$ws.Range("E37").Value = "  -5.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.677"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.102.23"
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01613"
$ws.Range("E40").Value = "  -4.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.899"
$ws.Range("E41").Value = "  -6.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8611"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.96"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.844.11"
$ws.Range("E45").Value = "  -5.29%  "
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.57"
$ws.Range("E47").Value = "  -5.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.026"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4308"
$ws.Range("E50").Value = "  -3.69%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05214"
$ws.Range("E51").Value = "  -4.21%  "
